# Updates odds values on Sheet1 as per the latest Betfair Back/Lay refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.73
$ws.Range("G2").Value = 1.76
$ws.Range("I2").Value = 4.8
$ws.Range("L2").Value = 1.24
$ws.Range("N2").Value = 6.2
$ws.Range("P2").Value = 2.8
$ws.Range("R2").Value = 1.75
$ws.Range("U2").Value = 2.68
$ws.Range("V2").Value = 1.26
$ws.Range("AN2").Value = 7.4

# Row 3
$ws.Range("F3").Value = 1.86
$ws.Range("G3").Value = 1.98
$ws.Range("H3").Value = 4
$ws.Range("J3").Value = 3.85
$ws.Range("K3").Value = 4.4
$ws.Range("L3").Value = 1.31
$ws.Range("O3").Value = 1.23
$ws.Range("P3").Value = 2.26
$ws.Range("T3").Value = 1.65
$ws.Range("W3").Value = 2.02
$ws.Range("AB3").Value = 14
$ws.Range("AN3").Value = 10.5

# Row 6
$ws.Range("F6").Value = 2.26
$ws.Range("Q6").Value = 2.06
$ws.Range("T6").Value = 1.82
$ws.Range("Y6").Value = 13.5
$ws.Range("AE6").Value = 46
$ws.Range("AL6").Value = 40

# Row 7
$ws.Range("F7").Value = 1.56
$ws.Range("G7").Value = 1.58
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 4.8
$ws.Range("R7").Value = 1.58
$ws.Range("U7").Value = 2.22
$ws.Range("Y7").Value = 28
$ws.Range("AA7").Value = 210
$ws.Range("AG7").Value = 9.800000000000001
$ws.Range("AN7").Value = 6.6

# Row 8
$ws.Range("P8").Value = 2.28

# Row 9
$ws.Range("H9").Value = 3.15
$ws.Range("I9").Value = 3.2
$ws.Range("P9").Value = 2.06
$ws.Range("Q9").Value = 1.91
$ws.Range("AO9").Value = 30

# Row 10
$ws.Range("O10").Value = 1.22
$ws.Range("P10").Value = 2.38
$ws.Range("Z10").Value = 9
$ws.Range("AE10").Value = 14.5
$ws.Range("AH10").Value = 24
$ws.Range("AN10").Value = 140
